# Generate Report for Handoff
# Replaces the e2e test file GUID (and its content-hash-derived xliff
# filenames) with a freshly generated one, and bumps the handoff /
# handback timestamps to reflect the new run.

$wb = $excel.ActiveWorkbook

$oldGuid = "ad3498be-d35c-4155-a580-44f8c186f0d6"
$newGuid = "d1c77145-b09c-44b4-80e2-68135d36d963"

$oldHash = "7a02f23b06fdd3cc5a2e2cb7cddb0e653e619b9e"
$newHash = "4c5163c3fae4d8e30d64e079ba3c0f7c5a6cb013"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07e3952282cb00ad503de1518cdfdcd1b86f733d/e2e/$oldGuid.md"

function Update-Hyperlink($range, $display) {
    # Re-creating (rather than mutating in place) is required: this host
    # appends a stray override entry instead of editing the existing
    # <hyperlink> when only TextToDisplay is assigned.
    $range.Hyperlinks.Delete()
    $range.Worksheet.Hyperlinks.Add($range, $hyperlinkAddress, "", "", $display) | Out-Null
    $range.Font.Underline = $true
    $range.Font.Color = 15570276
}

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
Update-Hyperlink $wsOverview.Range("B2") "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-16 20:53:52"

# ---- Sheet "zh-cn" ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
Update-Hyperlink $wsZh.Range("A2") "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-16 20:53:46"

# ---- Sheet "de-de" ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
Update-Hyperlink $wsDe.Range("A2") "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-16 20:53:52"
